$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 291
$ws.Range("B291").Value = 6941375
$ws.Range("F291").Value = "Al Hazm"
$ws.Range("G291").Value = "Al Fateh SC"
$ws.Range("H291").Value = 2
$ws.Range("I291").Value = 0
$ws.Range("K291").Value = 5
$ws.Range("L291").Value = 3.9
$ws.Range("N291").Value = 3.1
$ws.Range("R291").Value = 1.925
$ws.Range("S291").Value = 1.875
$ws.Range("U291").Value = 1.825
$ws.Range("V291").Value = 1.975
$ws.Range("W291").Value = 2.1
$ws.Range("Z291").Value = 0.925
$ws.Range("AB291").Value = -1
$ws.Range("AC291").Value = 0.9750000000000001

# Row 292
$ws.Range("B292").Value = 6940761
$ws.Range("F292").Value = "Damac FC"
$ws.Range("G292").Value = "Al Ittihad Jeddah"
$ws.Range("H292").Value = 3
$ws.Range("I292").Value = 1
$ws.Range("K292").Value = 4.5
$ws.Range("L292").Value = 4.2
$ws.Range("N292").Value = 2.9
$ws.Range("R292").Value = 1.825
$ws.Range("S292").Value = 1.975
$ws.Range("U292").Value = 2
$ws.Range("V292").Value = 1.8
$ws.Range("W292").Value = 1.9
$ws.Range("Z292").Value = 0.825
$ws.Range("AB292").Value = 1
$ws.Range("AC292").Value = -1

# Row 294
$ws.Range("B294").Value = 6941378
$ws.Range("F294").Value = "Al Taawon Buraidah"
$ws.Range("G294").Value = "Al Fayha"
$ws.Range("H294").Value = 4
$ws.Range("I294").Value = 1
$ws.Range("J294").Value = "H"
$ws.Range("K294").Value = 1.75
$ws.Range("L294").Value = 3.6
$ws.Range("M294").Value = 4.2
$ws.Range("N294").Value = 1.7
$ws.Range("O294").Value = 3.75
$ws.Range("P294").Value = 4.5
$ws.Range("Q294").Value = -0.75
$ws.Range("R294").Value = 1.95
$ws.Range("S294").Value = 1.85
$ws.Range("T294").Value = 2.5
$ws.Range("U294").Value = 1.8
$ws.Range("V294").Value = 2
$ws.Range("W294").Value = 0.7
$ws.Range("Y294").Value = -1
$ws.Range("Z294").Value = 0.95
$ws.Range("AB294").Value = 0.8
$ws.Range("AC294").Value = -1

# Row 295
$ws.Range("B295").Value = 6941379
$ws.Range("F295").Value = "Al Khaleej Saihat"
$ws.Range("G295").Value = "Abha"
$ws.Range("H295").Value = 3
$ws.Range("K295").Value = 2.15
$ws.Range("L295").Value = 3.5
$ws.Range("M295").Value = 3
$ws.Range("N295").Value = 1.75
$ws.Range("O295").Value = 3.8
$ws.Range("P295").Value = 4
$ws.Range("R295").Value = 1.975
$ws.Range("S295").Value = 1.825
$ws.Range("T295").Value = 3
$ws.Range("U295").Value = 2
$ws.Range("V295").Value = 1.8
$ws.Range("W295").Value = 0.75
$ws.Range("Z295").Value = 0.9750000000000001
$ws.Range("AB295").Value = 1

# Row 296
$ws.Range("B296").Value = 6940759
$ws.Range("F296").Value = "Al Taee"
$ws.Range("G296").Value = "Al Hilal Riyadh"
$ws.Range("H296").Value = 1
$ws.Range("I296").Value = 2
$ws.Range("J296").Value = "A"
$ws.Range("K296").Value = 9.5
$ws.Range("L296").Value = 6.5
$ws.Range("M296").Value = 1.2
$ws.Range("N296").Value = 15
$ws.Range("O296").Value = 9
$ws.Range("P296").Value = 1.111
$ws.Range("Q296").Value = 2.5
$ws.Range("R296").Value = 1.9
$ws.Range("S296").Value = 1.9
$ws.Range("T296").Value = 3.75
$ws.Range("U296").Value = 1.975
$ws.Range("V296").Value = 1.825
$ws.Range("W296").Value = -1
$ws.Range("Y296").Value = 0.111
$ws.Range("Z296").Value = 0.8999999999999999
$ws.Range("AB296").Value = -1
$ws.Range("AC296").Value = 0.825

# Row 303
$ws.Range("B303").Value = 6940763
$ws.Range("F303").Value = "Al Hilal Riyadh"
$ws.Range("G303").Value = "Al Wehda Mecca"
$ws.Range("H303").Value = 2
$ws.Range("I303").Value = 0
$ws.Range("J303").Value = "H"
$ws.Range("K303").Value = 1.2
$ws.Range("L303").Value = 7
$ws.Range("M303").Value = 9.5
$ws.Range("N303").Value = 1.222
$ws.Range("O303").Value = 7
$ws.Range("P303").Value = 8.5
$ws.Range("Q303").Value = -2
$ws.Range("R303").Value = 1.85
$ws.Range("S303").Value = 1.95
$ws.Range("T303").Value = 3.75
$ws.Range("U303").Value = 1.9
$ws.Range("V303").Value = 1.9
$ws.Range("W303").Value = 0.222
$ws.Range("X303").Value = -1
$ws.Range("AC303").Value = 0.8999999999999999

# Row 304
$ws.Range("B304").Value = 6940818
$ws.Range("F304").Value = "Al Fateh SC"
$ws.Range("G304").Value = "Al Shabab Riyadh"
$ws.Range("H304").Value = 1
$ws.Range("J304").Value = "D"
$ws.Range("K304").Value = 2.55
$ws.Range("M304").Value = 2.45
$ws.Range("N304").Value = 2.375
$ws.Range("O304").Value = 3.5
$ws.Range("P304").Value = 2.625
$ws.Range("R304").Value = 1.8
$ws.Range("S304").Value = 2
$ws.Range("T304").Value = 2.75
$ws.Range("U304").Value = 1.775
$ws.Range("V304").Value = 2.025
$ws.Range("X304").Value = 2.5
$ws.Range("Y304").Value = -1
$ws.Range("Z304").Value = 0
$ws.Range("AA304").Value = 0
$ws.Range("AC304").Value = 1.025

# Row 305
$ws.Range("B305").Value = 6941381
$ws.Range("F305").Value = "Al Raed"
$ws.Range("G305").Value = "Damac FC"
$ws.Range("H305").Value = 0
$ws.Range("I305").Value = 1
$ws.Range("J305").Value = "A"
$ws.Range("K305").Value = 2.6
$ws.Range("L305").Value = 3.4
$ws.Range("M305").Value = 2.4
$ws.Range("N305").Value = 2.75
$ws.Range("O305").Value = 3.25
$ws.Range("P305").Value = 2.375
$ws.Range("Q305").Value = 0
$ws.Range("R305").Value = 2.025
$ws.Range("S305").Value = 1.775
$ws.Range("T305").Value = 2.5
$ws.Range("U305").Value = 1.95
$ws.Range("V305").Value = 1.85
$ws.Range("W305").Value = -1
$ws.Range("Y305").Value = 1.375
$ws.Range("Z305").Value = -1
$ws.Range("AA305").Value = 0.7749999999999999
$ws.Range("AC305").Value = 0.8500000000000001

# Row 311
$ws.Range("B311").Value = 6940819
$ws.Range("F311").Value = "Al Akhdoud"
$ws.Range("G311").Value = "Al Shabab Riyadh"
$ws.Range("H311").Value = 1
$ws.Range("I311").Value = 0
$ws.Range("J311").Value = "H"
$ws.Range("K311").Value = 2.5
$ws.Range("L311").Value = 3.4
$ws.Range("M311").Value = 2.6
$ws.Range("N311").Value = 2.7
$ws.Range("O311").Value = 3.25
$ws.Range("P311").Value = 2.5
$ws.Range("Q311").Value = 0
$ws.Range("R311").Value = 1.95
$ws.Range("S311").Value = 1.85
$ws.Range("T311").Value = 2.5
$ws.Range("U311").Value = 2
$ws.Range("V311").Value = 1.8
$ws.Range("W311").Value = 1.7
$ws.Range("Y311").Value = -1
$ws.Range("Z311").Value = 0.95
$ws.Range("AA311").Value = -1
$ws.Range("AB311").Value = -1
$ws.Range("AC311").Value = 0.8

# Row 313
$ws.Range("B313").Value = 6941388
$ws.Range("F313").Value = "Al Hazm"
$ws.Range("G313").Value = "Al Ahli Jeddah"
$ws.Range("H313").Value = 0
$ws.Range("I313").Value = 4
$ws.Range("J313").Value = "A"
$ws.Range("K313").Value = 5
$ws.Range("L313").Value = 4.75
$ws.Range("M313").Value = 1.444
$ws.Range("N313").Value = 5.5
$ws.Range("O313").Value = 5
$ws.Range("P313").Value = 1.363
$ws.Range("Q313").Value = 1.5
$ws.Range("R313").Value = 1.8
$ws.Range("S313").Value = 2
$ws.Range("T313").Value = 3.25
$ws.Range("U313").Value = 2.025
$ws.Range("V313").Value = 1.775
$ws.Range("W313").Value = -1
$ws.Range("Y313").Value = 0.363
$ws.Range("Z313").Value = -1
$ws.Range("AA313").Value = 1
$ws.Range("AB313").Value = 1.025
$ws.Range("AC313").Value = -1

# Row 341
$ws.Range("R341").Value = 1.8
$ws.Range("S341").Value = 2
